# Cp questionnaire usecontext (#393) - "add experimental" commit
# Updates the Metadata sheet of the ValueSet workbook:
#   - Experimental value -> "true"
#   - Date value         -> "2025-01-28T15:58:19+00:00"
#   - Description value  -> "ValueSet regroupant des valuesets du NOS pour le
#                             code de la division territoriale"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")   # Metadata / Property-Value sheet

# Helper: writes $text into $range as a genuine text (shared-string) value.
# Assigning "true"/"false" (or other auto-recognised literals) straight to
# Value/Value2 makes Excel coerce the cell to a Boolean, which is not what we
# want here (the source file stores it as plain text). Routing the text
# through a formula ("=...") that evaluates to a string, then copying just
# the *value* back into the destination cell, keeps it a real text cell
# while preserving the destination cell's existing style/formatting.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("Z100")
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# Row 7: Experimental | true
Set-TextValue $ws.Range("B7") "true"

# Row 8: Date | 2025-01-28T15:58:19+00:00
$ws.Range("B8").Value2 = "2025-01-28T15:58:19+00:00"

# Row 12: Description | ValueSet regroupant des valuesets du NOS pour le code de la division territoriale
$ws.Range("B12").Value2 = "ValueSet regroupant des valuesets du NOS pour le code de la division territoriale"

$excel.CutCopyMode = 0
